$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving numeric-looking strings
# (e.g. '213.45') as literal text instead of letting Excel coerce them to
# numbers (which would also introduce float rounding noise), while leaving
# the cell's style pointing back at the default 'Normal' style afterwards.
function Set-TextValue([string]$addr, [string]$val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '26.636.03'
Set-TextValue "E2" '  +1.74%  '

# Row 3
Set-TextValue "D3" '1.628.12'
Set-TextValue "E3" '  +1.83%  '

# Row 4
Set-TextValue "E4" '  -0.11%  '

# Row 5
Set-TextValue "D5" '213.45'

# Row 6
Set-TextValue "E6" '  -0.12%  '

# Row 7
Set-TextValue "D7" '0.488'
Set-TextValue "E7" '  +0.85%  '

# Row 8
Set-TextValue "E8" '  +0.80%  '

# Row 9
Set-TextValue "D9" '0.0619'
Set-TextValue "E9" '  +1.07%  '

# Row 10
Set-TextValue "D10" '18.98'
Set-TextValue "E10" '  +4.63%  '

# Row 11
Set-TextValue "D11" '0.0834'
Set-TextValue "E11" '  +2.56%  '

# Row 12
Set-TextValue "B12" 'WrappedEther'
Set-TextValue "C12" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D12" '1.691.71'
Set-TextValue "E12" '  +5.84%  '

# Row 13
Set-TextValue "B13" 'WrappedliquidstakedEther2.0'
Set-TextValue "C13" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D13" '1.855.48'
Set-TextValue "E13" '  +1.80%  '

# Row 14
Set-TextValue "E14" '  +0.77%  '

# Row 15
Set-TextValue "E15" '  +2.35%  '

# Row 16
Set-TextValue "D16" '26.607.30'
Set-TextValue "E16" '  +1.59%  '

# Row 17
Set-TextValue "D17" '63.05'
Set-TextValue "E17" '  +2.93%  '

# Row 18
Set-TextValue "E18" '  +0.68%  '

# Row 19
Set-TextValue "E19" '  +0.00%  '

# Row 20
Set-TextValue "D20" '206.34'
Set-TextValue "E20" '  +2.28%  '

# Row 21
Set-TextValue "D21" '4.32'
Set-TextValue "E21" '  +0.96%  '

# Row 22
Set-TextValue "D22" '9.41'
Set-TextValue "E22" '  +1.85%  '

# Row 23
Set-TextValue "E23" '  +2.02%  '

# Row 24
Set-TextValue "E24" '  -1.52%  '

# Row 25
Set-TextValue "D25" '145.67'
Set-TextValue "E25" '  +1.16%  '

# Row 26
Set-TextValue "E26" '  -0.13%  '

# Row 27
Set-TextValue "E27" '  -1.18%  '

# Row 28
Set-TextValue "D28" '15.37'
Set-TextValue "E28" '  +1.36%  '

# Row 29
Set-TextValue "D29" '6.64'
Set-TextValue "E29" '  +1.46%  '

# Row 30
Set-TextValue "D30" '0.0522'
Set-TextValue "E30" '  +7.27%  '

# Row 31
Set-TextValue "E31" '  +0.57%  '

# Row 32
Set-TextValue "D32" '3.23'
Set-TextValue "E32" '  +1.57%  '

# Row 33
Set-TextValue "D33" '2.94'
Set-TextValue "E33" '  +0.82%  '

# Row 34
Set-TextValue "E34" '  +1.77%  '

# Row 35
Set-TextValue "E35" '  -0.48%  '

# Row 36
Set-TextValue "D36" '1.162.50'
Set-TextValue "E36" '  +1.01%  '

# Row 37
Set-TextValue "D37" '0.0164'
Set-TextValue "E37" '  +0.32%  '

# Row 38
Set-TextValue "D38" '0.809'
Set-TextValue "E38" '  +2.31%  '

# Row 39
Set-TextValue "E39" '  -0.08%  '

# Row 40
Set-TextValue "E40" '  -0.23%  '

# Row 41
Set-TextValue "D41" '0.501'
Set-TextValue "E41" '  +0.87%  '

# Row 42
Set-TextValue "D42" '5.39'
Set-TextValue "E42" '  +3.31%  '

# Row 43
Set-TextValue "D43" '0.787'
Set-TextValue "E43" '  +0.93%  '

# Row 44
Set-TextValue "D44" '1.763.84'
Set-TextValue "E44" '  +1.54%  '

# Row 45
Set-TextValue "D45" '92.13'
Set-TextValue "E45" '  +0.47%  '

# Row 46
Set-TextValue "E46" '  +1.88%  '

# Row 47
Set-TextValue "D47" '54.40'
Set-TextValue "E47" '  +0.83%  '

# Row 48
Set-TextValue "D48" '0.0₆0102'
Set-TextValue "E48" '  +4.67%  '

# Row 49
Set-TextValue "D49" '0.0512'
Set-TextValue "E49" '  +1.09%  '

# Row 50
Set-TextValue "B50" 'Mantle'
Set-TextValue "C50" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D50" '0.409'
Set-TextValue "E50" '  +0.57%  '

# Row 51
Set-TextValue "B51" 'EnergySwap'
Set-TextValue "C51" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D51" '7.53'
Set-TextValue "E51" '  +4.41%  '
